$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) stays the same text, just re-asserting values ---
$ws.Range("A1").Value = "Country(en)"
$ws.Range("B1").Value = "MobileCode(en)"
$ws.Range("C1").Value = "State(en)"
$ws.Range("D1").Value = "City(en)"
$ws.Range("E1").Value = "Country(ar)"
$ws.Range("F1").Value = "MobileCode(ar)"
$ws.Range("G1").Value = "State(ar)"
$ws.Range("H1").Value = "City(ar)"

# --- Row 2: Doha / Ad Dawhah (existing row, gains a State(en)/State(ar) pair) ---
$ws.Range("A2").Value = "Qatar"
$ws.Range("B2").Value = 974
$ws.Range("C2").Value = "Ad Dawhah"
$ws.Range("D2").Value = "Doha"
$ws.Range("E2").Value = "الهند"
$ws.Range("F2").Value = 974
$ws.Range("G2").Value = "الدوحة"
$ws.Range("H2").Value = "الدوحة"

# --- Row 3: new Al Rayyan row, replacing the old stray "United Arab Emirates" cell ---
# Clear the old special formatting (white 7pt Courier New) that lived on A3 so the
# row reverts to plain/default styling like the rest of the data rows.
$ws.Range("A3").Style = "Normal"

$ws.Range("A3").Value = "Qatar"
$ws.Range("B3").Value = 974
$ws.Range("C3").Value = "Al Rayyan"
$ws.Range("D3").Value = "Al Rayyan"
$ws.Range("E3").Value = "الهند"
$ws.Range("F3").Value = 974
$ws.Range("G3").Value = "الريان"
$ws.Range("H3").Value = "الريان"

# --- Number format for the mobile-code columns (B, F) on the new row ---
$ws.Range("B3").NumberFormat = "\+0"
$ws.Range("F3").NumberFormat = "\+0"

# --- Restore the active-cell selection recorded in the saved file ---
$ws.Range("F6").Select()
